$d = $word.ActiveDocument

# The target run lives in the final paragraph of the document body
# (just before the body's sectPr). Update its text and append a new,
# blank paragraph (same tab stops / run-formatting as the original)
# right after it, by including the paragraph mark ("^p") in the
# replacement text of a Find/Replace.
$old = "}}: PS     OM     EM;"
$new = "}}: PS     OM     EM; Community Control: PS    EM; County Jail: PS   EM;^p"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $lastPara.Range

$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "Replaced: $found"
Write-Output "Paragraph count now: $($d.Paragraphs.Count)"
